# Adds the "November 18 2023" meeting minutes entry (row 16) and fills in the
# missing "Ended" time (column D) for the "November 16 2023" entry (row 15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Give D15 (the previously-empty "Ended" cell for the Nov 16 row) the same
# time number-format as the rest of column D/C by copying C15's format.
$ws.Range("C15").Copy()
$ws.Range("D15").PasteSpecial(-4122)  # xlPasteFormats

# Duplicate the formatting of the whole Nov 16 row (A15:E15) down into the
# new Nov 18 row (A16:E16) before filling in its values.
$ws.Range("A15:E15").Copy()
$ws.Range("A16:E16").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the new values. The order below matches the order new entries were
# appended to the shared-string table by the original author.
$ws.Range("A16").Value = "November 18 2023"
$ws.Range("C16").Value = "2:00PM"
$ws.Range("D15").Value = "2:45PM"
$ws.Range("D16").Value = "2:15PM"
$ws.Range("E16").Value = "Update each other on what was done and what needs to be done going forwards"
$ws.Range("B16").Value = "Sedat, David, Madison, Sean"

# Match the author's final selection/view state.
$ws.Range("E16").Select() | Out-Null
